$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -1
$ws.Range("F4").Value = -2
$ws.Range("F5").Value = -9
$ws.Range("F7").Value = -2
$ws.Range("F9").Value = -5
$ws.Range("F10").Value = -3
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = -1
$ws.Range("F15").Value = -6
$ws.Range("F16").Value = -6
